$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text/label/link cells and non-ambiguous numeric-looking text (safe to assign directly)
$ws.Range("D2").Value = "46.665.67"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").Value = "2.581.17"
$ws.Range("E3").Value = "  +10.42%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("E7").Value = "  +6.20%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +12.07%  "
$ws.Range("E10").Value = "  +12.29%  "
$ws.Range("E11").Value = "  +5.25%  "
$ws.Range("E12").Value = "  +12.66%  "
$ws.Range("D13").Value = "2.973.09"
$ws.Range("E13").Value = "  +10.35%  "
$ws.Range("D15").Value = "2.596.52"
$ws.Range("E15").Value = "  +10.85%  "
$ws.Range("E16").Value = "  +12.12%  "
$ws.Range("E17").Value = "  +10.15%  "
$ws.Range("D18").Value = "46.918.18"
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("E19").Value = "  +9.17%  "
$ws.Range("D20").Value = "0.0₃0999"
$ws.Range("E20").Value = "  +3.88%  "
$ws.Range("E21").Value = "  +9.98%  "
$ws.Range("E22").Value = "  +5.22%  "
$ws.Range("E23").Value = "  +4.06%  "
$ws.Range("E24").Value = "  +6.44%  "
$ws.Range("E25").Value = "  +11.28%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("E27").Value = "  +23.58%  "
$ws.Range("E28").Value = "  +2.73%  "
$ws.Range("E29").Value = "  +8.27%  "
$ws.Range("E30").Value = "  +3.68%  "
$ws.Range("E31").Value = "  +5.51%  "
$ws.Range("B32").Value = "WEMIXToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("E32").Value = "  +5.52%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E33").Value = "  +8.86%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E34").Value = "  +9.48%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E35").Value = "  +21.90%  "
$ws.Range("E36").Value = "  +3.48%  "
$ws.Range("E37").Value = "  +9.14%  "
$ws.Range("E38").Value = "  +4.13%  "
$ws.Range("E39").Value = "  +9.21%  "
$ws.Range("E40").Value = "  +10.04%  "
$ws.Range("E41").Value = "  +9.69%  "
$ws.Range("E42").Value = "  +11.71%  "
$ws.Range("D43").Value = "2.016.43"
$ws.Range("E43").Value = "  +9.60%  "
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("B45").Value = "BitcoinSV"
$ws.Range("C45").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("E45").Value = "  +2.75%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E46").Value = "  +36.63%  "
$ws.Range("E47").Value = "  +2.33%  "
$ws.Range("E48").Value = "  +11.68%  "
$ws.Range("E49").Value = "  +7.78%  "
$ws.Range("D50").Value = "2.835.88"
$ws.Range("E50").Value = "  +10.48%  "
$ws.Range("E51").Value = "  +10.33%  "

# Cells whose text content looks like a plain number (e.g. "8.00", "0.998").
# Force text storage via a temporary text number format so Excel keeps the exact
# string (including trailing zeros) as text rather than converting it to a number,
# then reset the style back to Normal so no stray style index is left on the cell.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.601"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.571"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0832"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.105"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.899"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "15.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "41.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.43"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0843"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "149.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.121"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.24"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0328"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "107.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.200"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.84"
$ws.Range("D51").Style = "Normal"
